$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" text on Hoja1!A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$cellA1 = $ws1.Range("A1")
$text = $cellA1.Value2
$text = $text.Replace("1000 Bs = 3.14 = 11730.96 pesos", "1000 Bs = 3.07 = 11437.29 pesos")
$text = $text.Replace("11730.96 pesos = 3.12 = 970.56 Bs", "11437.29 pesos = 3.05 = 964.54 Bs")
$cellA1.Value = $text

# --- Update the tasas sheet numeric cells ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 325.69
$ws2.Range("O10").Value = 3725.01
$ws2.Range("N12").Value = 3753
$ws2.Range("O12").Value = 316.5
